# "both mega beakraft and AI edits"
# Update the Mega Beakraft base-stat block (I3:N3) and the derived-at-level-100
# stat row (I5:N5) with the corrected numbers, and tidy up the "mega <stat>"
# label row (I4:N4) into a single filled-right formula so it matches the other
# label rows (A4:F4, I6:N6, A10:F10, A16:F16) which are already filled that way.
#
# Everything downstream (O3, O5 sums; I7:N7 "max" stats; J8:N8 stat-mult
# ratios) is a formula and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mega Beakraft base stats (row 3: hp, atk, def, speed, spatk, spdef) ---
$ws.Range("I3").Value = 75
$ws.Range("J3").Value = 75
$ws.Range("K3").Value = 75
$ws.Range("L3").Value = 40
$ws.Range("M3").Value = 125
$ws.Range("N3").Value = 95

# --- Mega Beakraft stats at level 100, 0 IV / 0 EV (row 5) ---
$ws.Range("I5").Value = 75
$ws.Range("J5").Value = 105
$ws.Range("K5").Value = 75
$ws.Range("L5").Value = 70
$ws.Range("M5").Value = 135
$ws.Range("N5").Value = 125

# --- Re-fill the "mega <stat>" label row as a single shared formula (I4:N4) ---
$ws.Range("I4:N4").Formula = "=""mega "" & I$2"
